$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "period" text + formatting for rows 2-9 (column B).
# Old shared string "Jul 2022 - Jun 2023 data" becomes a new string
# (a literal tab, newline, then "Oct 2022-Sep 2023") and the cells pick
# up a new number format / wrap-text style.
$rng = $ws.Range("B2:B9")
$rng.Value = "`t`nOct 2022-Sep 2023"
$rng.NumberFormat = "@"
$rng.WrapText = $true

# Update the worksheet view: scroll so row 9 is at the top and select B2:B9.
$ws.Activate()
$ws.Range("B2:B9").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
